$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H98").Value = 361843.53
$ws.Range("I98").Value = 508639.28
$ws.Range("J98").Value = 3009.4443
$ws.Range("K98").Value = 508639.28
$ws.Range("L98").Value = 3009.4443
$ws.Range("M98").Value = -507141.28
$ws.Range("N98").Value = -6005.4443
$ws.Range("H122").Value = 361843.53
$ws.Range("I122").Value = 508639.28
$ws.Range("J122").Value = 3009.4443
$ws.Range("K122").Value = 1525917.84
$ws.Range("L122").Value = 9028.332900000001
$ws.Range("M122").Value = -1523467.84
$ws.Range("N122").Value = -13928.3329
$ws.Range("H129").Value = 1009.9677
$ws.Range("I129").Value = 339.5
$ws.Range("J129").Value = 1170.88
$ws.Range("K129").Value = 1018.5
$ws.Range("L129").Value = 3512.64
$ws.Range("M129").Value = 3981.5
$ws.Range("N129").Value = -13512.64
$ws.Range("H138").Value = 6581430
$ws.Range("I138").Value = 2475.348
$ws.Range("K138").Value = 7426.044
$ws.Range("M138").Value = -2286.044

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1131
$ws.Range("I122").Value = 1053.25
$ws.Range("J122").Value = 1397.5714
$ws.Range("K122").Value = 3159.75
$ws.Range("L122").Value = 4192.7142
$ws.Range("M122").Value = -709.75
$ws.Range("N122").Value = -9092.7142
$ws.Range("H132").Value = 1981.5
$ws.Range("I132").Value = 1329.4651
$ws.Range("J132").Value = 4530.364
$ws.Range("K132").Value = 3988.3953
$ws.Range("L132").Value = 13591.092
$ws.Range("M132").Value = -1458.3953
$ws.Range("N132").Value = -18651.092

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10195.889
$ws.Range("J86").Value = 10760.909
$ws.Range("L86").Value = 10760.909
$ws.Range("N86").Value = -13006.909
$ws.Range("H89").Value = 10195.889
$ws.Range("J89").Value = 10760.909
$ws.Range("L89").Value = 53804.545
$ws.Range("N89").Value = -65036.545
$ws.Range("H134").Value = 4542.8125
$ws.Range("I134").Value = 2040
$ws.Range("J134").Value = 5680.4546
$ws.Range("K134").Value = 6120
$ws.Range("L134").Value = 17041.3638
$ws.Range("M134").Value = -3585
$ws.Range("N134").Value = -22111.3638

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1206.2858
$ws.Range("I19").Value = 1206.2858
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1206.2858
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1036.2858
$ws.Range("N19").ClearContents()
$ws.Range("H24").Value = 1206.2858
$ws.Range("I24").Value = 1206.2858
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1206.2858
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -1036.2858
$ws.Range("N24").ClearContents()
$ws.Range("H132").Value = 2176.8157
$ws.Range("I132").Value = 1323.7037
$ws.Range("J132").Value = 4270.8184
$ws.Range("K132").Value = 3971.1111
$ws.Range("L132").Value = 12812.4552
$ws.Range("M132").Value = -1441.1111
$ws.Range("N132").Value = -17872.4552

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 50
$ws.Range("K26").Value = 150
$ws.Range("M26").Value = 138
$ws.Range("H33").Value = 234.21428
$ws.Range("I33").Value = 104.7
$ws.Range("J33").Value = 558
$ws.Range("K33").Value = 628.2
$ws.Range("L33").Value = 3348
$ws.Range("M33").Value = -345.2
$ws.Range("N33").Value = -3914
$ws.Range("H38").Value = 122.6129
$ws.Range("I38").Value = 163.46153
$ws.Range("J38").Value = 93.111115
$ws.Range("K38").Value = 490.38459
$ws.Range("L38").Value = 279.333345
$ws.Range("M38").Value = -143.38459
$ws.Range("N38").Value = -973.333345
$ws.Range("H41").Value = 614.2857
$ws.Range("J41").Value = 645.4545000000001
$ws.Range("L41").Value = 1936.3635
$ws.Range("N41").Value = -2612.3635
$ws.Range("H141").Value = 6134.5454
$ws.Range("I141").Value = 7170
$ws.Range("J141").Value = 5542.857
$ws.Range("K141").Value = 21510
$ws.Range("L141").Value = 16628.571
$ws.Range("M141").Value = -16330
$ws.Range("N141").Value = -26988.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1072.6595
$ws.Range("I102").Value = 960.7027
$ws.Range("J102").Value = 1486.9
$ws.Range("K102").Value = 960.7027
$ws.Range("L102").Value = 1486.9
$ws.Range("M102").Value = 661.2973
$ws.Range("N102").Value = -4730.9
$ws.Range("H113").Value = 1400.579
$ws.Range("I113").Value = 1198.875
$ws.Range("J113").Value = 1547.2727
$ws.Range("K113").Value = 1198.875
$ws.Range("L113").Value = 1547.2727
$ws.Range("M113").Value = 971.125
$ws.Range("N113").Value = -5887.2727
$ws.Range("H123").Value = 10302.167
$ws.Range("J123").Value = 10302.167
$ws.Range("L123").Value = 10302.167
$ws.Range("N123").Value = -15202.167

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3281
$ws.Range("I7").Value = 2480.2
$ws.Range("J7").Value = 3531.25
$ws.Range("K7").Value = 2480.2
$ws.Range("L7").Value = 3531.25
$ws.Range("M7").Value = -2368.2
$ws.Range("N7").Value = -3755.25
$ws.Range("H22").Value = 706.26666
$ws.Range("I22").Value = 844.44446
$ws.Range("J22").Value = 499
$ws.Range("K22").Value = 844.44446
$ws.Range("L22").Value = 499
$ws.Range("M22").Value = -549.44446
$ws.Range("N22").Value = -1089
$ws.Range("H27").Value = 706.26666
$ws.Range("I27").Value = 844.44446
$ws.Range("J27").Value = 499
$ws.Range("K27").Value = 844.44446
$ws.Range("L27").Value = 499
$ws.Range("M27").Value = -737.44446
$ws.Range("N27").Value = -713
$ws.Range("H122").Value = 3024.8057
$ws.Range("I122").Value = 2220.9285
$ws.Range("J122").Value = 3536.3635
$ws.Range("K122").Value = 6662.7855
$ws.Range("L122").Value = 10609.0905
$ws.Range("M122").Value = -4212.7855
$ws.Range("N122").Value = -15509.0905
$ws.Range("H126").Value = 3281
$ws.Range("I126").Value = 2480.2
$ws.Range("J126").Value = 3531.25
$ws.Range("K126").Value = 7440.599999999999
$ws.Range("L126").Value = 10593.75
$ws.Range("M126").Value = -4970.599999999999
$ws.Range("N126").Value = -15533.75
$ws.Range("H132").Value = 2890.689
$ws.Range("I132").Value = 2087.9429
$ws.Range("J132").Value = 5700.3
$ws.Range("K132").Value = 6263.8287
$ws.Range("L132").Value = 17100.9
$ws.Range("M132").Value = -3733.8287
$ws.Range("N132").Value = -22160.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 33290
$ws.Range("J125").Value = 33290
$ws.Range("L125").Value = 33290
$ws.Range("N125").Value = -43130
$ws.Range("H126").Value = 79029.234
$ws.Range("I126").Value = 113153.336
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 339460.008
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -336990.008
$ws.Range("N126").Value = -11690
$ws.Range("H132").Value = 12822740
$ws.Range("I132").Value = 18520236
$ws.Range("J132").Value = 3371.6667
$ws.Range("K132").Value = 55560708
$ws.Range("L132").Value = 10115.0001
$ws.Range("M132").Value = -55558178
$ws.Range("N132").Value = -15175.0001
